$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old row 107 ("Fuyu" / Primera) so the
# existing rows 107-109 shift down to 109-111, and the two freshly
# inserted rows (107-108) can hold the new "Mankaki" (2022-06-12) data.
$ws.Rows.Item(107).Insert()
$ws.Rows.Item(107).Insert()

# New row 107: Mankaki / Primera
$ws.Range("A107").Value = 8
$ws.Range("B107").Value = "Terminal La Palmera de La Serena"
$ws.Range("C107").Value = "Coquimbo"
$ws.Range("D107").Value = 44714
$ws.Range("E107").Value = 4
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100107
$ws.Range("H107").Value = "Otros"
$ws.Range("I107").Value = 100107001
$ws.Range("J107").Value = "Caqui"
$ws.Range("K107").Value = "Mankaki"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 16
$ws.Range("N107").Value = 330000
$ws.Range("O107").Value = 340000
$ws.Range("P107").Value = 335000
$ws.Range("Q107").Value = "`$/bins (450 kilos)"
$ws.Range("R107").Value = "Región de O'Higgins"
$ws.Range("S107").Value = 744
$ws.Range("T107").Value = 450

# New row 108: Mankaki / Segunda
$ws.Range("A108").Value = 8
$ws.Range("B108").Value = "Terminal La Palmera de La Serena"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44714
$ws.Range("E108").Value = 4
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100107
$ws.Range("H108").Value = "Otros"
$ws.Range("I108").Value = 100107001
$ws.Range("J108").Value = "Caqui"
$ws.Range("K108").Value = "Mankaki"
$ws.Range("L108").Value = "Segunda"
$ws.Range("M108").Value = 16
$ws.Range("N108").Value = 290000
$ws.Range("O108").Value = 300000
$ws.Range("P108").Value = 295000
$ws.Range("Q108").Value = "`$/bins (450 kilos)"
$ws.Range("R108").Value = "Región de O'Higgins"
$ws.Range("S108").Value = 656
$ws.Range("T108").Value = 450
